$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($i=1; $i -le 8; $i++) {
  $c = $ws.Columns.Item($i)
  Write-Host "col $i ColumnWidth: $($c.ColumnWidth) Width: $($c.Width)"
}
